$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.239.08"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "1.644.88"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.33"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.875.74"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.659.40"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.42"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "27.234.86"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.69"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("E22").Value = "  +7.40%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.64"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.52"
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").Value = "1.276.58"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.545"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +6.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.31"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "1.785.62"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.98"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.96"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("E51").Value = "  +0.77%  "
